$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '29.371.51'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').Value = '1.600.31'
$ws.Range('E3').Value = '  +2.28%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  +0.06%  '
Set-TextValue 'D5' '211.97'
$ws.Range('E5').Value = '  +0.64%  '
Set-TextValue 'D6' '0.517'
$ws.Range('E6').Value = '  +6.35%  '
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '26.41'
$ws.Range('E8').Value = '  +5.55%  '
Set-TextValue 'D9' '43.42'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('E10').Value = '  +2.30%  '
Set-TextValue 'D11' '0.0598'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('E12').Value = '  +1.63%  '
$ws.Range('D13').Value = '1.832.74'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').Value = '1.621.29'
$ws.Range('E14').Value = '  +3.92%  '
$ws.Range('D15').Value = '29.445.66'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('E16').Value = '  +3.59%  '
$ws.Range('E17').Value = '  +1.37%  '
Set-TextValue 'D18' '63.16'
$ws.Range('E18').Value = '  +2.91%  '
Set-TextValue 'D19' '239.23'
$ws.Range('E19').Value = '  +3.40%  '
Set-TextValue 'D20' '7.61'
$ws.Range('E20').Value = '  +3.37%  '
$ws.Range('D21').Value = '0.0₃0687'
$ws.Range('E21').Value = '  +1.93%  '
Set-TextValue 'D22' '0.999'
$ws.Range('E22').Value = '  -0.05%  '
$ws.Range('E23').Value = '  +1.61%  '
Set-TextValue 'D24' '9.11'
$ws.Range('E24').Value = '  +1.69%  '
Set-TextValue 'D25' '2.08'
$ws.Range('E25').Value = '  -1.59%  '
Set-TextValue 'D26' '154.59'
$ws.Range('E26').Value = '  +2.91%  '
Set-TextValue 'D27' '15.23'
$ws.Range('E27').Value = '  +3.07%  '
$ws.Range('E28').Value = '  +4.41%  '
Set-TextValue 'D29' '6.35'
$ws.Range('E29').Value = '  +2.17%  '
Set-TextValue 'D30' '1.00'
$ws.Range('E30').Value = '  -0.02%  '
Set-TextValue 'D31' '0.0471'
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('E32').Value = '  -0.16%  '
$ws.Range('E33').Value = '  +1.39%  '
Set-TextValue 'D34' '3.09'
$ws.Range('E34').Value = '  +3.69%  '
$ws.Range('D35').Value = '1.409.66'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').Value = '  -0.86%  '
$ws.Range('E37').Value = '  +2.62%  '
Set-TextValue 'D38' '2.81'
$ws.Range('E38').Value = '  +5.30%  '
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('E40').Value = '  +2.30%  '
Set-TextValue 'D41' '0.534'
$ws.Range('E41').Value = '  +3.38%  '
Set-TextValue 'D42' '1.97'
$ws.Range('E42').Value = '  +1.03%  '
Set-TextValue 'D43' '0.0485'
$ws.Range('E43').Value = '  +5.61%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D44' '0.999'
$ws.Range('E44').Value = '  -0.08%  '
$ws.Range('B45').Value = 'ARBITRUM'
$ws.Range('C45').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D45' '0.791'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D46' '51.98'
$ws.Range('E46').Value = '  +20.46%  '
Set-TextValue 'D47' '65.44'
$ws.Range('E47').Value = '  +2.62%  '
Set-TextValue 'D48' '5.26'
$ws.Range('E48').Value = '  -0.10%  '
$ws.Range('D49').Value = '1.742.51'
$ws.Range('E49').Value = '  +2.51%  '
$ws.Range('B50').Value = 'WEMIXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D50' '0.854'
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('B51').Value = 'Quant'
$ws.Range('C51').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D51' '86.47'
$ws.Range('E51').Value = '  +1.29%  '
